$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 579
$ws.Range("F3").Value = 5474
$ws.Range("F4").Value = 61
$ws.Range("F5").Value = 457
$ws.Range("F8").Value = 380
$ws.Range("F9").Value = 1352
$ws.Range("F11").Value = 1885
$ws.Range("F12").Value = 3074
$ws.Range("F13").Value = 1910
$ws.Range("F15").Value = 57
$ws.Range("F16").Value = 191
$ws.Range("F17").Value = 29
$ws.Range("F18").Value = 142
$ws.Range("F19").Value = 650
$ws.Range("F20").Value = 973
$ws.Range("F21").Value = 350
$ws.Range("F22").Value = 48
$ws.Range("F23").Value = 3543
$ws.Range("F24").Value = 1113
$ws.Range("F25").Value = 2808
$ws.Range("F26").Value = 281
$ws.Range("F27").Value = 1972
$ws.Range("F28").Value = 4061
$ws.Range("F30").Value = 915
$ws.Range("F31").Value = 465
$ws.Range("F32").Value = 1293
$ws.Range("F33").Value = 61
$ws.Range("F36").Value = 1267
$ws.Range("F37").Value = 62
$ws.Range("F38").Value = 1041
$ws.Range("F39").Value = 672
$ws.Range("F40").Value = 536
$ws.Range("F41").Value = 409
$ws.Range("F42").Value = 25
$ws.Range("F43").Value = 312
$ws.Range("F44").Value = 3564

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 21
$ws.Range("F10").Value = 905
$ws.Range("F11").Value = 32
$ws.Range("F22").Value = 37
$ws.Range("F23").Value = 29
$ws.Range("F24").Value = 4

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 579
$ws.Range("F3").Value = 579
$ws.Range("F4").Value = 5474
$ws.Range("F5").Value = 61
$ws.Range("F6").Value = 21
$ws.Range("F9").Value = 380
$ws.Range("F10").Value = 1352
$ws.Range("F11").Value = 3074
$ws.Range("F13").Value = 1910
$ws.Range("F15").Value = 57
$ws.Range("F17").Value = 191
$ws.Range("F18").Value = 905
$ws.Range("F19").Value = 32
$ws.Range("F21").Value = 142
$ws.Range("F22").Value = 973
$ws.Range("F23").Value = 350
$ws.Range("F24").Value = 3543
$ws.Range("F27").Value = 1114
$ws.Range("F28").Value = 2809
$ws.Range("F29").Value = 1972
$ws.Range("F30").Value = 4061
$ws.Range("F33").Value = 915
$ws.Range("F34").Value = 1293
$ws.Range("F38").Value = 1267
$ws.Range("F39").Value = 62
$ws.Range("F40").Value = 1041
$ws.Range("F42").Value = 672
$ws.Range("F44").Value = 409
$ws.Range("F45").Value = 37
$ws.Range("F46").Value = 29
$ws.Range("F48").Value = 312
$ws.Range("F49").Value = 3564
